# The deck's "datetimeFigureOut" date placeholder (Insert > Header & Footer
# style date field) is cached as "9/29/14" on the slide master and on every
# one of its slide layouts. Bump the cached date to "9/30/14" everywhere it
# appears, the way PowerPoint itself would when the date field is edited.

$p = $ppt.ActivePresentation
$newDate = "9/30/14"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($container, $text) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $phType = -1
            try {
                $phType = $shp.PlaceholderFormat.Type
            } catch {
                $phType = -1
            }
            if ($phType -eq $ppPlaceholderDate) {
                $shp.TextFrame.TextRange.Text = $text
            }
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master $newDate

# Every slide layout owned by the master has its own copy of the placeholder.
$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DatePlaceholder $layouts.Item($j) $newDate
}
